
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13 (shifts rows 13-23 down to 14-24)
$ws.Rows.Item(13).Insert()

# Fix row 10 (A10 "Objetivos:") - B10/C10 currently hold docentes data, replace with PT objectives text
$ws.Cells.Item(10,2).Value = '1. Aproximar o aluno do seu aparelho vocal, ao nível da expressão falada ou cantada, em sua expressão individual ou coletiva (coral). 2. Propiciar ao aluno o contato com o repertório coral e seu papel no desenvolvimento da linguagem musical. 3. Propiciar ao aluno a oportunidade de fazer música vocal em conjunto, com o aprendizado técnico de parâmetros como: afinação, precisão, equilíbrio, fraseado etc.'
$ws.Cells.Item(10,3).Value = '1. Aproximar o aluno do seu aparelho vocal, ao nível da expressão falada ou cantada, em sua expressão individual ou coletiva (coral). 2. Propiciar ao aluno o contato com o repertório coral e seu papel no desenvolvimento da linguagem musical. 3. Propiciar ao aluno a oportunidade de fazer música vocal em conjunto, com o aprendizado técnico de parâmetros como: afinação, precisão, equilíbrio, fraseado etc.'

# New row 13: docentes responsaveis data in B13/C13 (A13 left blank)
$ws.Cells.Item(13,2).Value = '8452037 - Elisabeth Pinheiro da Silva Kondracki de Alcantara'
$ws.Cells.Item(13,3).Value = '8452037 - Elisabeth Pinheiro da Silva Kondracki de Alcantara'
$ws.Cells.Item(13,1).Clear()
$ws.Cells.Item(13,2).WrapText = $true
$ws.Cells.Item(13,2).VerticalAlignment = -4160
$ws.Cells.Item(13,2).Font.Bold = $false

# Fix row 14 (A14 "Programa resumido:") - B14/C14 currently hold 'Semestral', replace with PT short syllabus
$ws.Cells.Item(14,2).Value = 'Retomada e aperfeiçoamento dos repertórios e pressupostos técnicos da disciplina precedente de Canto Coral. Classificação Vocal. Respiração para o canto. Colocação da emissão em "Bocca Chiusa". O canto coral em uníssono.O canto coral em cânone. O canto coral em outras formações polifônicas.Leitura coral.'
$ws.Cells.Item(14,3).Value = 'Retomada e aperfeiçoamento dos repertórios e pressupostos técnicos da disciplina precedente de Canto Coral. Classificação Vocal. Respiração para o canto. Colocação da emissão em "Bocca Chiusa". O canto coral em uníssono.O canto coral em cânone. O canto coral em outras formações polifônicas.Leitura coral.'

# Fix row 16 (A16 "Programa:") - B16/C16 currently hold a date, replace with PT long syllabus
$ws.Cells.Item(16,2).Value = 'Retomada e aperfeiçoamento dos repertórios e pressupostos técnicos da disciplina precedente de Canto Coral. Classificação Vocal. Respiração para o canto - exercícios para localização da respiração baixa e média. Apoio e coluna de ar. Colocação da emissão em "Bocca Chiusa" relaxamento da mandíbula, posicionamento da língua, suspensão do palato mole, exercícios de percepção do local onde a voz está se colocando, conexão do apoio e emissão vocal, passagem da "Bocca Chiusa" para vogais e outros sons nasais e guturais. O canto coral em uníssono: afinação, uniformidade tímbrica, precisão rítmica. O canto coral em cânone. O canto coral em outras formações polifônicas. Leitura coral: testagem da escuta harmônica e afinação à 1ª vista a várias vozes, memorização, fundamentos de teoria musical. Montagem e aperfeiçoamento de peças musicais - promovendo a aplicação das técnicas aprendidas. Conexão entre diafragma e emissão vocal.'
$ws.Cells.Item(16,3).Value = 'Retomada e aperfeiçoamento dos repertórios e pressupostos técnicos da disciplina precedente de Canto Coral. Classificação Vocal. Respiração para o canto - exercícios para localização da respiração baixa e média. Apoio e coluna de ar. Colocação da emissão em "Bocca Chiusa" relaxamento da mandíbula, posicionamento da língua, suspensão do palato mole, exercícios de percepção do local onde a voz está se colocando, conexão do apoio e emissão vocal, passagem da "Bocca Chiusa" para vogais e outros sons nasais e guturais. O canto coral em uníssono: afinação, uniformidade tímbrica, precisão rítmica. O canto coral em cânone. O canto coral em outras formações polifônicas. Leitura coral: testagem da escuta harmônica e afinação à 1ª vista a várias vozes, memorização, fundamentos de teoria musical. Montagem e aperfeiçoamento de peças musicais - promovendo a aplicação das técnicas aprendidas. Conexão entre diafragma e emissão vocal.'

# Fix row 19 (A19 "Método:") - B19/C19 currently hold docentes data, replace with the "A cada semestre..." text
$ws.Cells.Item(19,2).Value = 'A cada semestre é proposto um programa com cerca de 8 (oito) peças, sendo duas ou três de semestres anteriores e, consequentemente, cinco ou seis inéditas – a ser apresentado pelo CORAL da EEL-USP em performances públicas definidas durante o período letivo.'
$ws.Cells.Item(19,3).Value = 'A cada semestre é proposto um programa com cerca de 8 (oito) peças, sendo duas ou três de semestres anteriores e, consequentemente, cinco ou seis inéditas – a ser apresentado pelo CORAL da EEL-USP em performances públicas definidas durante o período letivo.'

# Fix row 20 (A20 "Critério:") - B20/C20 currently hold "A cada semestre..." text, replace with "Sendo uma atividade..." text
$ws.Cells.Item(20,2).Value = 'Sendo uma atividade prática e de grupo, fica inviável a realização de provas ou outras formas similares de avaliação. Esta se dará no dia a dia do aluno, levando em conta: assiduidade, pontualidade e material completo na pasta; participação construtiva em sala de aula e nas apresentações públicas - prontidão, envolvimento e seu real aproveitamento vocal e musical.'
$ws.Cells.Item(20,3).Value = 'Sendo uma atividade prática e de grupo, fica inviável a realização de provas ou outras formas similares de avaliação. Esta se dará no dia a dia do aluno, levando em conta: assiduidade, pontualidade e material completo na pasta; participação construtiva em sala de aula e nas apresentações públicas - prontidão, envolvimento e seu real aproveitamento vocal e musical.'

# Fix row 21 (A21 "Norma de recuperacao:") - B21/C21 currently hold "Sendo uma atividade..." text, replace with "Nao tem"
$ws.Cells.Item(21,2).Value = 'Não tem'
$ws.Cells.Item(21,3).Value = 'Não tem'

# Fix row 22 (A22 "Bibliografia:") - B22/C22 currently hold "Nao tem", replace with the full bibliography text
$ws.Cells.Item(22,2).Value = 'BEHLAU, Mara; REHDER, Maria Inês. Higiene vocal para o Canto Coral. Rio de Janeiro: Revinter, 1997.COELHO, Helena Wöhl. Técnica vocal para coros. 7ª Edição. São Leopoldo: Sinodal, 1994.CORBIN, Lynn Ann. Vocal pedagogy in the choral rehearsal: The selected concepts on choral tone quality, understanding of the singing process, and attitudes toward choir participation. 1982. 119 f. Tese (Doutorado em Filosofia) - Graduate School of The Ohio State University, Ohio State University, 1982.COSTA, Paulo Rubens Moraes. Diagnose em Canto Coral: parâmetros para análise e ferramentas para a avaliação. 2005. Dissertação (Mestrado em Musicologia) - Departamento de Música, Escola de Comunicações e Artes, Universidade de São Paulo, São Paulo, 2005.DRAHAN, Snizhana. Ouvir a voz: a percepção da produção vocal pelo regente coral - método eformação. 2007. 146 f. Dissertação (Mestrado em Musicologia) - Departamento de Música, Escola de Comunicações e Artes, Universidade de São Paulo, São Paulo, 2007.FERNANDES, Angelo José. O regente moderno e a construção da sonoridade coral: uma metodologia de preparo vocal para coros. 2009. 475 f. Tese (Doutorado em Música) - Instituto de Artes, Universidade Federal de Campinas, Campinas, 2009.HERR, Martha. Considerações para a classificação da voz do coralista. In: FERREIRA, LésliePiccolotto et al. Voz profissional: o profissional da voz. Carapicuíba: Pró-fono DepartamentoEditorial, 1995. p. 51-56.MILLER, Richard. The structure of singing. Boston, Massachusetts: Schirmer, 1986.RAMOS, Marco Antonio da Silva. Canto Coral: do repertório temático à construção do programa. 1988. 492f. Dissertação (Mestrado em Artes) - Departamento de Biblioteconomia e Documentação, Escola de Comunicações e Artes, Universidade de São Paulo, São Paulo, 1988.______. Memorial. 2011. Memorial (Professor titular) - Departamento de Música, Escola de Comunicações e Artes, Universidade de São Paulo, São Paulo, 2011.______. O ensino da regência coral. 2003. 107f. Tese (Livre Docência) - Departamento de Música, Escola de Comunicações e Artes, Universidade de São Paulo, São Paulo, 2003.VIDEIRA JR., Mário Rodrigues. Educação musical através do coro: A experiência no Projeto Comunicantus da ECA - USP. 2001. 65 f. Trabalho de Conclusão de Curso (Licenciatura em Música) - Departamento de Música, Escola de Comunicações e Artes, Universidade de São Paulo, São Paulo, 2001.'
$ws.Cells.Item(22,3).Value = 'BEHLAU, Mara; REHDER, Maria Inês. Higiene vocal para o Canto Coral. Rio de Janeiro: Revinter, 1997.COELHO, Helena Wöhl. Técnica vocal para coros. 7ª Edição. São Leopoldo: Sinodal, 1994.CORBIN, Lynn Ann. Vocal pedagogy in the choral rehearsal: The selected concepts on choral tone quality, understanding of the singing process, and attitudes toward choir participation. 1982. 119 f. Tese (Doutorado em Filosofia) - Graduate School of The Ohio State University, Ohio State University, 1982.COSTA, Paulo Rubens Moraes. Diagnose em Canto Coral: parâmetros para análise e ferramentas para a avaliação. 2005. Dissertação (Mestrado em Musicologia) - Departamento de Música, Escola de Comunicações e Artes, Universidade de São Paulo, São Paulo, 2005.DRAHAN, Snizhana. Ouvir a voz: a percepção da produção vocal pelo regente coral - método eformação. 2007. 146 f. Dissertação (Mestrado em Musicologia) - Departamento de Música, Escola de Comunicações e Artes, Universidade de São Paulo, São Paulo, 2007.FERNANDES, Angelo José. O regente moderno e a construção da sonoridade coral: uma metodologia de preparo vocal para coros. 2009. 475 f. Tese (Doutorado em Música) - Instituto de Artes, Universidade Federal de Campinas, Campinas, 2009.HERR, Martha. Considerações para a classificação da voz do coralista. In: FERREIRA, LésliePiccolotto et al. Voz profissional: o profissional da voz. Carapicuíba: Pró-fono DepartamentoEditorial, 1995. p. 51-56.MILLER, Richard. The structure of singing. Boston, Massachusetts: Schirmer, 1986.RAMOS, Marco Antonio da Silva. Canto Coral: do repertório temático à construção do programa. 1988. 492f. Dissertação (Mestrado em Artes) - Departamento de Biblioteconomia e Documentação, Escola de Comunicações e Artes, Universidade de São Paulo, São Paulo, 1988.______. Memorial. 2011. Memorial (Professor titular) - Departamento de Música, Escola de Comunicações e Artes, Universidade de São Paulo, São Paulo, 2011.______. O ensino da regência coral. 2003. 107f. Tese (Livre Docência) - Departamento de Música, Escola de Comunicações e Artes, Universidade de São Paulo, São Paulo, 2003.VIDEIRA JR., Mário Rodrigues. Educação musical através do coro: A experiência no Projeto Comunicantus da ECA - USP. 2001. 65 f. Trabalho de Conclusão de Curso (Licenciatura em Música) - Departamento de Música, Escola de Comunicações e Artes, Universidade de São Paulo, São Paulo, 2001.'
